$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# --- New row 32: Credits ---
$ws.Range("A32").Value = "Credits"
$ws.Range("B32").Value = "Credits"
$ws.Range("C32").Value = "Mitwirkende"
$ws.Range("D32").Value = "Załoga"

# --- New row 33: HowToPlay ---
$ws.Range("A33").Value = "HowToPlay"
$ws.Range("B33").Value = "How to Play"
$ws.Range("C33").Value = "Spielanleitung"
$ws.Range("D33").Value = "Jak Grać"

# --- New row 34: OpenSourceNotice ---
$ws.Range("A34").Value = "OpenSourceNotice"
$ws.Range("B34").Value = 'This game is open source. You can find it on <color=#A52A2A><link source="githuburl">GitHub  \uf35d</link></color>'
$ws.Range("C34").Formula = '=""'
$ws.Range("D34").Formula = '=""'

# --- Existing placeholder cells that now get real values ---
$ws.Range("D8").Value = "Start"
$ws.Range("D7").Value = "Pestępój"

# --- Move active selection ---
[void]$ws.Range("D11").Select()
